$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = ""

$ws.Range("H69").Value = 8503.75
$ws.Range("J69").Value = 9671.666999999999
$ws.Range("L69").Value = 29015.001
$ws.Range("N69").Value = -30763.001

$ws.Range("H72").Value = 8503.75
$ws.Range("J72").Value = 9671.666999999999
$ws.Range("L72").Value = 87045.003
$ws.Range("N72").Value = -95781.003

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""

$ws.Range("H86").Value = 446.66666
$ws.Range("I86").Value = 446.66666
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 446.66666
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 676.33334
$ws.Range("N86").Value = ""

$ws.Range("H89").Value = 446.66666
$ws.Range("I89").Value = 446.66666
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 2233.3333
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 3382.6667
$ws.Range("N89").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17999
$ws.Range("J61").Value = 17999
$ws.Range("L61").Value = 17999
$ws.Range("N61").Value = -18423

$ws.Range("H74").Value = 958.8
$ws.Range("I74").Value = 958.8
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 958.8
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -84.79999999999995
$ws.Range("N74").Value = ""

$ws.Range("H77").Value = 958.8
$ws.Range("I77").Value = 958.8
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4794
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -426
$ws.Range("N77").Value = ""

$ws.Range("H92").Value = 49774
$ws.Range("J92").Value = 49774
$ws.Range("L92").Value = 49774
$ws.Range("N92").Value = -54766

$ws.Range("H132").Value = 5607.6
$ws.Range("I132").Value = 5607.6
$ws.Range("K132").Value = 16822.8
$ws.Range("M132").Value = -14292.8

$ws.Range("H136").Value = 17999
$ws.Range("J136").Value = 17999
$ws.Range("L136").Value = 53997
$ws.Range("N136").Value = -59097

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 49800
$ws.Range("J68").Value = 49800
$ws.Range("L68").Value = 49800
$ws.Range("N68").Value = -51298

$ws.Range("H71").Value = 49800
$ws.Range("J71").Value = 49800
$ws.Range("L71").Value = 149400
$ws.Range("N71").Value = -156888

$ws.Range("H107").Value = 1390
$ws.Range("I107").Value = 1390
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1390
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 530
$ws.Range("N107").Value = ""

$ws.Range("H132").Value = 2855.7778
$ws.Range("I132").Value = 2855.7778
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8567.3334
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6037.3334
$ws.Range("N132").Value = ""

$ws.Range("H134").Value = 5333
$ws.Range("I134").Value = 10000
$ws.Range("J134").Value = 2999.5
$ws.Range("K134").Value = 30000
$ws.Range("L134").Value = 8998.5
$ws.Range("M134").Value = -27465
$ws.Range("N134").Value = -14068.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 60
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""

$ws.Range("H12").Value = 45.53846
$ws.Range("I12").Value = 52.2
$ws.Range("J12").Value = 41.375
$ws.Range("K12").Value = 156.6
$ws.Range("L12").Value = 124.125
$ws.Range("M12").Value = 16.39999999999998
$ws.Range("N12").Value = -470.125

$ws.Range("H34").Value = 1750
$ws.Range("J34").Value = 1750
$ws.Range("L34").Value = 5250
$ws.Range("N34").Value = -5418

$ws.Range("H39").Value = 4999.6665
$ws.Range("J39").Value = 4999.6665
$ws.Range("L39").Value = 14998.9995
$ws.Range("N39").Value = -15586.9995

$ws.Range("H131").Value = 823.3333
$ws.Range("I131").Value = 823.3333
$ws.Range("K131").Value = 2469.9999
$ws.Range("M131").Value = 2570.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3250
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -5996

$ws.Range("H83").Value = 3250
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -29984

$ws.Range("H126").Value = 2074.5
$ws.Range("J126").Value = 3150
$ws.Range("L126").Value = 9450
$ws.Range("N126").Value = -14390

$ws.Range("H132").Value = 4933.3335
$ws.Range("I132").Value = 4933.3335
$ws.Range("K132").Value = 14800.0005
$ws.Range("M132").Value = -12270.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 940.2
$ws.Range("I22").Value = 875
$ws.Range("K22").Value = 875
$ws.Range("M22").Value = -580

$ws.Range("H27").Value = 940.2
$ws.Range("I27").Value = 875
$ws.Range("K27").Value = 875
$ws.Range("M27").Value = -768

$ws.Range("H46").Value = 8599.556
$ws.Range("I46").Value = 9333.166999999999
$ws.Range("J46").Value = 7132.3335
$ws.Range("K46").Value = 9333.166999999999
$ws.Range("L46").Value = 7132.3335
$ws.Range("M46").Value = -9145.166999999999
$ws.Range("N46").Value = -7508.3335

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = ""

$ws.Range("H69").Value = 17165.166
$ws.Range("J69").Value = 17165.166
$ws.Range("L69").Value = 17165.166
$ws.Range("N69").Value = -18663.166

$ws.Range("H72").Value = 17165.166
$ws.Range("J72").Value = 17165.166
$ws.Range("L72").Value = 51495.49800000001
$ws.Range("N72").Value = -58983.49800000001

$ws.Range("H81").Value = 492.5
$ws.Range("I81").Value = 492.5
$ws.Range("K81").Value = 985
$ws.Range("M81").Value = 76

$ws.Range("H84").Value = 492.5
$ws.Range("I84").Value = 492.5
$ws.Range("K84").Value = 4925
$ws.Range("M84").Value = 379

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""

$ws.Range("H122").Value = 1640.6
$ws.Range("I122").Value = 1640.6
$ws.Range("K122").Value = 4921.799999999999
$ws.Range("M122").Value = -2471.799999999999
